$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (No 17)
$ws.Range("C18").Value = 17
$ws.Range("D18").Value = 'LIBRO BLANCO de los herbolarios y las plantas medicinales'
$ws.Range("E18").Value = 2007
$ws.Range("F18").Value = 'Fundación Salud y Naturaleza (S.N.)'
$ws.Range("H18").Value = 'España'
$ws.Range("I18").Value = "No"
$ws.Range("K18").Value = 'https://www.fitoterapia.net/archivos/200701/260307libro-2.pdf?1'

# Row 19 (No 18)
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 'EL GRAN LIBRO DE LA MEDICINA CHINA'
$ws.Range("E19").Value = 2003
$ws.Range("F19").Value = 'Ediciones URANO, S.A'
$ws.Range("H19").Value = 'España'
$ws.Range("I19").Value = "Si"
$ws.Range("K19").Value = 'http://bibliosjd.org/wp-content/uploads/2017/03/El-Gran-Libro-De-La-Medicina-China.pdf'
$ws.Range("G19").Font.Name = "SabonLTStd-Roman"
$ws.Range("G19").Font.Color = 0

# Row 20 (No 19)
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 'EL GRAN LIBRO DE LA MEDICINA NATURAL'
$ws.Range("F20").Value = 'Ediciones Masters'
$ws.Range("G20").Font.Name = "Arial"
$ws.Range("G20").Font.Size = 14
$ws.Range("G20").Font.Color = 0
$ws.Rows.Item(20).RowHeight = 17.4

# Row 21 (No 20)
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 'Manual para la práctica de la Medicina Natural y Tradicional'
$ws.Range("E21").Value = 2014
$ws.Range("F21").Value = 'Editorial Ciencias Médicas'
$ws.Range("I21").Value = "Si"
$ws.Range("K21").Value = 'https://instituciones.sld.cu/fcmdoct/files/2019/10/manual_medtrad_completo.pdf'

# Row 22 (No 21)
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 'Plantas medicinales de La Matamba y El Piñonal, municipio de Jamapa, Veracruz'
$ws.Range("E22").Value = 2015
$ws.Range("F22").Value = 'Instituto de Ecología A. C. (INECOL)'
$ws.Range("G22").Value = '978-607-7579-44-1'
$ws.Range("I22").Value = "No"

# Row 23 (No 22)
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 'Plantas medicinales del Patio de Ángel albino corzo, chiaPas'
$ws.Range("E23").Value = 2020
$ws.Range("F23").Value = 'Universidad Autónoma de Chiapas'
$ws.Range("H23").Value = '978-607-561-075-7'
$ws.Range("I23").Value = "No"
$ws.Range("K23").Value = 'https://editorial.unach.mx/documentos/digitales/_libs/plantasmedicinales.pdf'
$ws.Hyperlinks.Add($ws.Range("K23"), 'https://editorial.unach.mx/documentos/digitales/_libs/plantasmedicinales.pdf')

# Row 24 (No 23)
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 'Contribución al conocimiento de las plantas medicinales del municipio de Tlatlauquitepec, Puebla'
$ws.Range("E24").Value = 2022
$ws.Range("F24").Value = 'Ciencia Latina Revista Científica Multidisciplinar'
$ws.Hyperlinks.Add($ws.Range("G24"), "http://dx.doi.org/10.37811/cl_rcm.v6i6.3741", "", "", "http://dx.doi.org/10.37811/cl_rcm.v6i6.3741")
$ws.Range("G24").Value = '10.37811/cl_rcm.v6i6.3741'
$ws.Range("H24").Value = 'Mexico'
$ws.Range("I24").Value = 'CC BY 4.0'
$ws.Range("K24").Value = 'https://www.researchgate.net/publication/366057090_Contribucion_al_conocimiento_de_las_plantas_medicinales_del_municipio_de_Tlatlauquitepec_Puebla/fulltext/638fe03f484e65005bec8985/Contribucion-al-conocimiento-de-las-plantas-medicinales-del-municipio-de-Tlatlauquitepec-Puebla.pdf?origin=publicationDetail&_sg%5B0%5D=jYRMsuvIrjqqIV7VxQZEe3Obc5lx6Wcjrg3H8wvQFc2L_h0oo1zJbHtAc67eIeE4-iOTkDSu5TJvuic_vLPYlg.yDBk41gY1IZGKVAVqqLEPnUKDcmcuieL78i_TRGXP80VwQHRFc-MmrCBUVWd8rrMrW90JOgQ97H1YH_w2XOP4Q&_sg%5B1%5D=9UzfNT5FSWL9OpTEs5-JhJDupExQZDXmluX6173mQ-EvFoqvZ0_WfpxoYiD9-Z1mFjzIT2-_QmFnY2YHqWqCLhGeOXnrvIvmeal5h3AnsOlg.yDBk41gY1IZGKVAVqqLEPnUKDcmcuieL78i_TRGXP80VwQHRFc-MmrCBUVWd8rrMrW90JOgQ97H1YH_w2XOP4Q&_iepl=&_rtd=eyJjb250ZW50SW50ZW50IjoibWFpbkl0ZW0ifQ%3D%3D&_tp=eyJjb250ZXh0Ijp7ImZpcnN0UGFnZSI6InB1YmxpY2F0aW9uIiwicGFnZSI6InB1YmxpY2F0aW9uIiwicG9zaXRpb24iOiJwYWdlSGVhZGVyIn19'

# Selection update
$ws.Range("D20").Select()
